# Fix formatting of floating point numbers that were scrapped with
# Spanish/Argentine style separators (e.g. "46.974,45") so that they use a
# plain decimal-point style (e.g. "46974.45"): drop the "." thousands
# separator and turn the "," decimal separator into a ".".
#
# Also fix two "Razon social" entries that used a comma to separate two
# contractor names, which was ambiguous with the decimal separator; those
# commas become periods as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the "Importe" column keeps storing plain text (these values are
# not real numbers in the workbook - they are text such as "380,00") instead
# of letting Excel re-interpret the new, unambiguous-looking numeric strings
# (e.g. "380.00") as actual numbers.
$importeRange = $ws.Range("H2:H74")
$importeRange.NumberFormat = "@"

for ($row = 2; $row -le 74; $row++) {
    $cell = $ws.Cells.Item($row, 8)
    $text = [string]$cell.Value()
    $fixed = $text.Replace(".", "").Replace(",", ".")
    $cell.Value = $fixed
}

# Replace the ambiguous commas in these two contractor-name cells with
# periods.
$ws.Range("E36").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E45").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E37").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
